$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto market data refresh (prices & 1h volume %) + two name/link row swaps
$updates = @(
    @{ Row = 2; D = "65.318.32"; E = "  -0.07%  " },
    @{ Row = 3; D = "3.533.78"; E = "  +2.96%  " },
    @{ Row = 4; D = "1.00"; E = "  +0.11%  " },
    @{ Row = 5; D = "595.50"; E = "  +0.19%  " },
    @{ Row = 6; D = "139.02"; E = "  +2.89%  " },
    @{ Row = 7; D = "3.533.61"; E = "  +2.96%  " },
    @{ Row = 8; E = "  +0.31%  " },
    @{ Row = 9; D = "0.496"; E = "  +1.25%  " },
    @{ Row = 10; D = "0.125"; E = "  +3.86%  " },
    @{ Row = 11; D = "7.23"; E = "  -3.04%  " },
    @{ Row = 12; D = "0.392"; E = "  +3.93%  " },
    @{ Row = 13; D = "4.141.18"; E = "  +3.42%  " },
    @{ Row = 14; D = "0.0000187"; E = "  +4.12%  " },
    @{ Row = 15; D = "27.03"; E = "  +2.11%  " },
    @{ Row = 16; D = "3.523.42"; E = "  +2.34%  " },
    @{ Row = 17; E = "  +1.35%  " },
    @{ Row = 18; D = "65.214.14"; E = "  +0.01%  " },
    @{ Row = 19; D = "10.23"; E = "  +1.98%  " },
    @{ Row = 20; D = "5.86"; E = "  +2.54%  " },
    @{ Row = 21; D = "14.24"; E = "  +3.97%  " },
    @{ Row = 22; D = "395.64"; E = "  +0.79%  " },
    @{ Row = 23; D = "0.571"; E = "  +4.91%  " },
    @{ Row = 24; D = "74.52"; E = "  +1.89%  " },
    @{ Row = 25; D = "3.683.53"; E = "  +3.25%  " },
    @{ Row = 26; D = "1.00"; E = "  +0.03%  " },
    @{ Row = 27; D = "0.0000114"; E = "  +9.58%  " },
    @{ Row = 28; D = "7.76"; E = "  +8.18%  " },
    @{ Row = 29; D = "0.998"; E = "  +0.10%  " },
    @{ Row = 30; D = "2.26"; E = "  +1.01%  " },
    @{ Row = 31; D = "8.28"; E = "  +1.69%  " },
    @{ Row = 32; D = "3.556.85"; E = "  +3.67%  " },
    @{ Row = 33; E = "  +0.03%  " },
    @{ Row = 34; D = "23.80"; E = "  +5.64%  " },
    @{ Row = 35; D = "0.145"; E = "  +1.11%  " },
    @{ Row = 36; D = "1.23"; E = "  -0.02%  " },
    @{ Row = 37; B = "Aptos"; C = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D = "6.99"; E = "  +2.48%  " },
    @{ Row = 38; B = "Monero"; C = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D = "170.21"; E = "  -1.56%  " },
    @{ Row = 39; D = "1.54"; E = "  +1.11%  " },
    @{ Row = 40; D = "4.91"; E = "  +2.48%  " },
    @{ Row = 41; D = "0.0800"; E = "  +4.08%  " },
    @{ Row = 42; B = "EnergySwap"; C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D = "26.98"; E = "  +23.23%  " },
    @{ Row = 43; B = "Mantle"; C = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"; D = "0.820"; E = "  +1.01%  " },
    @{ Row = 44; D = "42.83"; E = "  -2.10%  " },
    @{ Row = 45; D = "1.00"; E = "  +0.20%  " },
    @{ Row = 46; D = "4.43"; E = "  +0.63%  " },
    @{ Row = 47; E = "  +10.69%  " },
    @{ Row = 48; D = "1.67"; E = "  +3.96%  " },
    @{ Row = 49; D = "6.82"; E = "  +4.20%  " },
    @{ Row = 50; D = "2.340.03"; E = "  +6.62%  " },
    @{ Row = 51; D = "2.14"; E = "  +0.43%  " },
)

# Ensure the Price/Volume columns are written back as plain text (matching the
# original inline-string cells) rather than being auto-coerced to numbers.
$numRange = $ws.Range("D2:E51")
$numRange.NumberFormat = "@"

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("B")) {
        $ws.Range("B$r").Value = $u.B
    }
    if ($u.ContainsKey("C")) {
        $ws.Range("C$r").Value = $u.C
    }
    if ($u.ContainsKey("D")) {
        $ws.Range("D$r").Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $ws.Range("E$r").Value = $u.E
    }
}

# Restore the default cell style so no extra text-format styling lingers.
$numRange.Style = "Normal"
